$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "B11" "63.96"
Set-TextValue "C11" "5.34"

Set-TextValue "B12" "24.88"
Set-TextValue "C12" "41.74"
Set-TextValue "D12" "66.62"

Set-TextValue "B14" "87.05"
Set-TextValue "C14" "12.58"
Set-TextValue "D14" "99.63"

Set-TextValue "B37" "11.78"
Set-TextValue "C37" "5.12"
Set-TextValue "D37" "16.89"

Set-TextValue "B38" "16.14"
Set-TextValue "C38" "50.51"
Set-TextValue "D38" "66.65"

Set-TextValue "B44" "11.32"
Set-TextValue "C44" "47.55"
Set-TextValue "D44" "58.88"
